# Implement WaveEquation.integrate output for Dirichlet boundary conditions:
# insert a new header row at the top of the (x, y) data table and label the
# two columns "x" and "y". All existing data rows shift down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 1, pushing all data down.
$ws.Rows("1:1").Insert()

# Label the new header row.
$ws.Range("A1").Value = "x"
$ws.Range("B1").Value = "y"
